# Insert a new weekly record for Poroto verde (Comercializadora del Agro de
# Limarí) at row 144. Excel's Insert() shifts the existing rows 144:211 down
# to 145:212, preserving all of their data/styles automatically (it matches
# the "shift rows down, fill row 144 with new data" pattern described by the
# diff, which culminates with the former last row (211) being duplicated as
# the new last row (212)).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 144..211 down to 145..212.
$ws.Rows(144).Insert()

# Populate the newly freed row 144 with the new observation.
$ws.Cells.Item(144, 1).Value = 2
$ws.Cells.Item(144, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(144, 3).Value = "Coquimbo"
$ws.Cells.Item(144, 4).Value = 44846
$ws.Cells.Item(144, 5).Value = 4
$ws.Cells.Item(144, 6).Value = 100112031
$ws.Cells.Item(144, 7).Value = "Poroto verde"
$ws.Cells.Item(144, 8).Value = "Magnum"
$ws.Cells.Item(144, 9).Value = "Primera"
$ws.Cells.Item(144, 10).Value = 700
$ws.Cells.Item(144, 11).Value = 37000
$ws.Cells.Item(144, 12).Value = 40000
$ws.Cells.Item(144, 13).Value = 38500
$ws.Cells.Item(144, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(144, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(144, 16).Value = 1540
$ws.Cells.Item(144, 17).Value = 25
$ws.Cells.Item(144, 18).Value = "Hortaliza"
